# Highlight the best (lowest) fit-error value among columns C/D/E for each
# data row with a new light-blue fill, matching the "Init repo for sim on 47"
# commit. The previously orange "Total" rows keep their orange fill, and the
# plain data cells that lose their now-redundant General-number-format style
# fall back to the workbook's default (Normal) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Light-blue "Lighter 60%" tint of the theme's Accent1 color (5B9BD5 -> BDD7EE),
# expressed as the packed BGR long Excel's Interior.Color expects
# (R=0xBD=189, G=0xD7=215, B=0xEE=238 -> 189 + 215*256 + 238*65536).
$highlightColor = 15652797

# Row-by-row winner (minimum value among C/D/E) that gets the new highlight.
$bestFitCells = @("D4", "D5", "D6", "D7", "D8", "E9", "E10", "E11", "E12")
foreach ($addr in $bestFitCells) {
    $ws.Range($addr).Interior.Color = $highlightColor
}

# The remaining C/D/E data cells in rows 4-9 no longer need their old
# (visually no-op) explicit-number-format style; reset them to the default
# workbook style.
$plainCells = @("C4", "E4", "C5", "E5", "C6", "E6", "C7", "E7", "C8", "E8", "C9", "D9")
foreach ($addr in $plainCells) {
    $ws.Range($addr).Style = "Normal"
}

# Selection at save time.
[void]$ws.Range("N19").Select()
